# Apply "email sent on 06022020" update:
# Fills Sheet1 rows 2-8 with parsed-name / username / domain / group formulas
# driven off new Email (col D) and Password (col J) values, and adds the
# "Capgemini" lookup value on Sheet2!C2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Populate raw data (order matters: it drives shared-string insertion order) ---
$ws1.Range("D2").Value = "sourabh.awasthi@capgemini.com"
$ws1.Range("J2").Value = "[Any2m(J43F*"

$ws1.Range("D3").Value = "sandipan.deb@capgemini.com"
$ws1.Range("J3").Value = 'iptAP7Y$OEx{'

$ws2.Range("C2").Value = "Capgemini"

$ws1.Range("D4").Value = "biswajit.deb@capgemini.com"
$ws1.Range("J4").Value = 'bZv$tr486biN'

$ws1.Range("D5").Value = "debanjan.das@capgemini.com"
$ws1.Range("J5").Value = 'wE?}?"5+y6tE'

$ws1.Range("D6").Value = "dhiraj.kajari@capgemini.com"
$ws1.Range("J6").Value = 'VpYIVt=nI@$v'

$ws1.Range("D7").Value = "mayur.bhorkar@capgemini.com"
$ws1.Range("J7").Value = 't@c%**O*T@Jo'

$ws1.Range("D8").Value = "manoj-kumar.b.s@capgemini.com"
$ws1.Range("J8").Value = 'xU")eOY4[N#3'

# --- K / M plain-value columns ---
$ws1.Range("K2:K8").Value = 80
$ws1.Range("M2:M8").Value = $true

# --- Shared formulas across rows 2:7 ---
$ws1.Range("A2:A7").Formula = "=PROPER(IFERROR(LEFT(C2,FIND(CHAR(46),C2)-1),C2))"
$ws1.Range("B2:B7").Formula = '=IFERROR(PROPER(RIGHT(C2,LEN(C2)-FIND("@",SUBSTITUTE(C2,".","@",((LEN(C2)-LEN(SUBSTITUTE(C2,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C2:C7").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D2,FIND(CHAR(64),D2)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws1.Range("E2:E7").Formula = "=LEFT(H2,FIND(CHAR(46),H2)-1)"
$ws1.Range("F2:F7").Formula = '=CONCATENATE("ITPartner\",I2)'
$ws1.Range("H2:H7").Formula = "=RIGHT(D2,LEN(D2)-FIND(CHAR(64),D2))"
$ws1.Range("I2:I7").Formula = "=PROPER(E2)"
$ws1.Range("P2:P7").Formula = "=COUNTIF(D:D,D2)"

# --- Row 8 gets its own (non-shared) copies of the same formulas ---
$ws1.Range("A8").Formula = "=PROPER(IFERROR(LEFT(C8,FIND(CHAR(46),C8)-1),C8))"
$ws1.Range("B8").Formula = '=IFERROR(PROPER(RIGHT(C8,LEN(C8)-FIND("@",SUBSTITUTE(C8,".","@",((LEN(C8)-LEN(SUBSTITUTE(C8,".","")))/LEN("\")))))), "Unknown")'
$ws1.Range("C8").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D8,FIND(CHAR(64),D8)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws1.Range("E8").Formula = "=LEFT(H8,FIND(CHAR(46),H8)-1)"
$ws1.Range("F8").Formula = '=CONCATENATE("ITPartner\",I8)'
$ws1.Range("H8").Formula = "=RIGHT(D8,LEN(D8)-FIND(CHAR(64),D8))"
$ws1.Range("I8").Formula = "=PROPER(E8)"
$ws1.Range("P8").Formula = "=COUNTIF(D:D,D8)"

$excel.Calculate()
